# Insert a new weekly price record as row 228 in the "Vega Monumental
# Concepción - Tomate" sheet. Existing rows 228-233 shift down to 229-234.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by one row, opening up row 228 for the
# new record (mirrors Excel's "Insert Sheet Rows").
$ws.Rows.Item(228).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(228, 1).Value  = 11
$ws.Cells.Item(228, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(228, 3).Value  = "Bíobío"
$ws.Cells.Item(228, 4).Value  = 44448
$ws.Cells.Item(228, 5).Value  = 8
$ws.Cells.Item(228, 6).Value  = 100112020
$ws.Cells.Item(228, 7).Value  = "Tomate"
$ws.Cells.Item(228, 8).Value  = "Larga vida"
$ws.Cells.Item(228, 9).Value  = "Segunda"
$ws.Cells.Item(228, 10).Value = 300
$ws.Cells.Item(228, 11).Value = 18000
$ws.Cells.Item(228, 12).Value = 18000
$ws.Cells.Item(228, 13).Value = 18000
$ws.Cells.Item(228, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(228, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(228, 16).Value = 1000
$ws.Cells.Item(228, 17).Value = 18
$ws.Cells.Item(228, 18).Value = "Hortaliza"
